# Applies the edits described by the diff:
#  1. Adds a new "2020" data column (column N) with values mirroring the
#     existing 2010-2019 columns, copying formatting from column M.
#  2. Leaves the sheet's active selection on M25 (as recorded by Excel in
#     the saved sheetView).
#  3. Corrects a typo in the recorded absolute save path (best effort -
#     this is an Excel-maintained MRU-style field, not normally exposed
#     through the object model, so this line is a no-op on hosts that
#     don't surface it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "2020" column (column N) ---------------------------------------
$ws.Range("N3").Value = 2020

$ws.Range("N5").Value  = 2198.6999999999998
$ws.Range("N6").Value  = 132.69999999999999
$ws.Range("N7").Value  = 242.9
$ws.Range("N8").Value  = 203.3
$ws.Range("N9").Value  = 202.8
$ws.Range("N10").Value = 284.7
$ws.Range("N11").Value = 294.89999999999998
$ws.Range("N12").Value = 802.5
$ws.Range("N13").Value = 28.1
$ws.Range("N14").Value = 6.8

$ws.Range("N16").Value = 27.4
$ws.Range("N17").Value = 17.5
$ws.Range("N18").Value = 24.7
$ws.Range("N19").Value = 31.5
$ws.Range("N20").Value = 30.4
$ws.Range("N21").Value = 24.8
$ws.Range("N22").Value = 30.7
$ws.Range("N23").Value = 30.1
$ws.Range("N24").Value = 21.2
$ws.Range("N25").Value = 11.6

# Row 15 is a spacer row with no value in column N either, but (like the
# rest of that row) the cell still carries formatting - the PasteSpecial
# loop below creates it (empty, formatted) by copying format from M15.

# Mirror the formatting (number format, borders, etc.) of the matching
# column-M cell onto every new column-N cell we just populated.
$formatRows = 3,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25
foreach ($r in $formatRows) {
    $ws.Range("M$r").Copy()
    $ws.Range("N$r").PasteSpecial(-4122) # xlPasteFormats
}
$excel.CutCopyMode = 0

# --- Sheet selection ------------------------------------------------------
$ws.Range("M25").Select()

# --- Cosmetic absPath fix (best effort; see note above) -------------------
$wb.AbsPath = "C:\Users\korozbaeva\Desktop\Показатели ЦУР для Платформы\Национальные показатели ЦУР\"
